# Odoo pharma-pos demo data: drop the previous sample rows and replace them
# with a single demo record ("Aciclovir" / "200 mg Tablet"), matching the
# commit "removed previous demo data and added one record for each model".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 becomes the one remaining demo product (Name / Size columns);
# the numeric Store Code / TGP Code / Is Sold columns (C2:E2 = 1) are
# already correct and are left untouched.
$ws.Range("A2").Value = "Aciclovir"
$ws.Range("B2").Value = "200 mg Tablet"

# Rows 3-6 held the rest of the old demo catalogue (Colchisin, Losartan, ...).
# Clear them out entirely so only the header + the single new record remain.
$ws.Range("A3:E6").ClearContents()

# Leave the selection where the author left it when they saved the file.
$ws.Range("B8").Select()
